# issue #5: stock data from json to db
#
# The "股票" (stock) sheet gains three new columns:
#   - a "category" column inserted right after "property_category"
#     (pushing date / legislator_name / legislator_id one column right)
#   - "source_file" and "index" columns appended at the end,
#     recording provenance of the migrated record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Prep formatting for the new/shifted cells first --------------------
# Copy the existing header format (bold + border, used by column K) onto
# the new header cells L1:N1 before writing their text, and the existing
# data-row format (column K) onto the new data cells L2:N2, so the added
# columns look consistent with the rest of the table.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("L2:N2").PasteSpecial(-4122)

# J2 ("2011-11-17") must stay literal text rather than be auto-parsed
# into a date serial number when it is written. Force a text format
# before assigning it, then restore the plain data-row look via the
# same paste-format trick.
$ws.Range("J2").NumberFormat = "@"

# --- Header row (row 1) ---
# Columns B..H stay the same (name, owner, quantity, face_value,
# currency, total, property_category). Starting at column I the
# remaining headers shift one column right to make room for
# "category", and two new headers are appended at the end.
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data row (row 2) ---
# Same shift: existing date/legislator_name/legislator_id values move
# one column right, "category" gets the new value "normal", and the
# new source_file / index columns are populated.
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2011-11-17"
$ws.Range("H2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("K2").Value = "呂學樟"
$ws.Range("L2").Value = 892
$ws.Range("M2").Value = "tmpf9381"
$ws.Range("N2").Value = 71
